# Update "想去人数" (number of people interested) figures across sheets
# as part of refreshing the scraped data (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 60
$ws1.Range("F3").Value = 122
$ws1.Range("F4").Value = 1958
$ws1.Range("F5").Value = 312
$ws1.Range("F6").Value = 56
$ws1.Range("F8").Value = 2043
$ws1.Range("F9").Value = 10360
$ws1.Range("F11").Value = 148
$ws1.Range("F15").Value = 7210
$ws1.Range("F16").Value = 1107
$ws1.Range("F17").Value = 681
$ws1.Range("F18").Value = 122
$ws1.Range("F20").Value = 261

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 17

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 60
$ws4.Range("F3").Value = 122
$ws4.Range("F4").Value = 1958
$ws4.Range("F5").Value = 312
$ws4.Range("F6").Value = 56
$ws4.Range("F7").Value = 17
$ws4.Range("F9").Value = 2043
$ws4.Range("F12").Value = 10360
$ws4.Range("F14").Value = 148
$ws4.Range("F18").Value = 7210
$ws4.Range("F19").Value = 1107
$ws4.Range("F20").Value = 681
$ws4.Range("F21").Value = 122
$ws4.Range("F23").Value = 261
